# prescriptions.xlsx update:
#  - row 2 (E2/F2/G2) gets new BEGIN_TIME/END_TIME/DOSE_MG values
#  - three more BloodPressureReading-style rows (3,4,5) are appended,
#    reusing the formatting already present in row 2
#  - selection moves to the new last cell, G5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- update existing row 2 ---
$ws.Range("E2").Value = 45583
$ws.Range("F2").Value = 45589.999305555553
$ws.Range("G2").Value = 100

# --- new row 3 ---
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 45590
$ws.Range("F3").Value = 45597.5
$ws.Range("G3").Value = 200

# --- new row 4 ---
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 45597.541666666664
$ws.Range("F4").Value = 45613.999305555553
$ws.Range("G4").Value = 150

# --- new row 5 (no END_TIME, i.e. F5 left blank) ---
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 45621
$ws.Range("G5").Value = 175

# Carry the row-2 cell formatting down into rows 3-5 so the new rows
# look like the existing ones (same numeric/date formats).
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3:A5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B2:D2").Copy() | Out-Null
$ws.Range("B3:D5").PasteSpecial(-4122) | Out-Null

$ws.Range("E2:F2").Copy() | Out-Null
$ws.Range("E3:F5").PasteSpecial(-4122) | Out-Null

$ws.Range("G2").Copy() | Out-Null
$ws.Range("G3:G5").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Move the selection to the new bottom-right data cell.
$ws.Range("G5").Select() | Out-Null
